$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Simon"
$ws.Range("B5").Value = "Setup project infrastructure "
$ws.Range("C5").Value = "Created code folder and structur. Added small demo script."
$ws.Range("D5").Value = 17.3
$ws.Range("E5").Value = "1h"

$ws.Range("E5").Select()
